$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row data per hunk 0 (old_start 2595)
$ws.Range("H40").Value = 7344.3
$ws.Range("J40").Value = 8689.200000000001
$ws.Range("L40").Value = 8689.200000000001
$ws.Range("N40").Value = -9039.200000000001

# Row data per hunk 1 (old_start 5228)
$ws.Range("H93").Value = 36000
$ws.Range("J93").Value = 36000
$ws.Range("L93").Value = 36000
$ws.Range("N93").Value = -40992

# Row data per hunk 2 (old_start 5424)
$ws.Range("H97").Value = 834.8333
$ws.Range("J97").Value = 834.8333
$ws.Range("L97").Value = 2504.4999
$ws.Range("N97").Value = -3496.4999

# Row data per hunk 3 (old_start 5577)
$ws.Range("H100").Value = 83.333336
$ws.Range("I100").Value = 83.333336
$ws.Range("K100").Value = 83.333336
$ws.Range("M100").Value = 457.666664

# Row data per hunk 4 (old_start 7381)
$ws.Range("H137").Value = 706387.4399999999
$ws.Range("I137").Value = 1585288.6
$ws.Range("K137").Value = 4755865.800000001
$ws.Range("M137").Value = -4753315.800000001

# Row data per hunk 5 (old_start 7577)
$ws.Range("H141").Value = 4834
$ws.Range("I141").Value = 3753.75
$ws.Range("J141").Value = 6994.5
$ws.Range("K141").Value = 11261.25
$ws.Range("L141").Value = 20983.5
$ws.Range("M141").Value = -6081.25
$ws.Range("N141").Value = -31343.5

$ws = $wb.Worksheets.Item("ARM")
# Row data per hunk 6 (old_start 7873)
$ws.Range("H5").Value = 351.83334
$ws.Range("I5").Value = 305
$ws.Range("K5").Value = 305
$ws.Range("M5").Value = -193

# Row data per hunk 7 (old_start 8709)
$ws.Range("H22").Value = 6167.9165
$ws.Range("I22").Value = 6167.9165
$ws.Range("K22").Value = 6167.9165
$ws.Range("M22").Value = -5868.9165

# Row data per hunk 8 (old_start 9208)
$ws.Range("H32").Value = 3706149.8
$ws.Range("I32").Value = 274.77274
$ws.Range("K32").Value = 274.77274
$ws.Range("M32").Value = 12.22726

# Row data per hunk 9 (old_start 9511)
$ws.Range("H38").Value = 3435.6667
$ws.Range("I38").Value = 3122.8
$ws.Range("J38").Value = 5000
$ws.Range("K38").Value = 3122.8
$ws.Range("L38").Value = 5000
$ws.Range("M38").Value = -2655.8
$ws.Range("N38").Value = -5934

# Row data per hunk 10 (old_start 9664)
$ws.Range("H41").Value = 1431.2
$ws.Range("I41").Value = 664
$ws.Range("J41").Value = 4500
$ws.Range("K41").Value = 664
$ws.Range("L41").Value = 4500
$ws.Range("M41").Value = -250
$ws.Range("N41").Value = -5328

# Row data per hunk 11 (old_start 11925)
$ws.Range("H88").Value = 962.3333
$ws.Range("I88").Value = 475
$ws.Range("K88").Value = 475
$ws.Range("M88").Value = -69

# Row data per hunk 12 (old_start 12069)
$ws.Range("H91").Value = 962.3333
$ws.Range("I91").Value = 475
$ws.Range("K91").Value = 475
$ws.Range("M91").Value = 929

# Row data per hunk 13 (old_start 12614)
$ws.Range("H102").Value = 6544.2856
$ws.Range("I102").Value = 905
$ws.Range("K102").Value = 905
$ws.Range("M102").Value = 717

# Row data per hunk 14 (old_start 13567)
$ws.Range("H122").Value = 1832.1333
$ws.Range("I122").Value = 1551.7273
$ws.Range("J122").Value = 2603.25
$ws.Range("K122").Value = 4655.1819
$ws.Range("L122").Value = 7809.75
$ws.Range("M122").Value = -2205.1819
$ws.Range("N122").Value = -12709.75

$ws = $wb.Worksheets.Item("BSM")
# Row data per hunk 15 (old_start 14727)
$ws.Range("H4").Value = 351.83334
$ws.Range("I4").Value = 305
$ws.Range("K4").Value = 305
$ws.Range("M4").Value = -190

# Row data per hunk 16 (old_start 19080)
$ws.Range("H94").Value = 592.6667
$ws.Range("I94").Value = 611.125
$ws.Range("J94").Value = 445
$ws.Range("K94").Value = 611.125
$ws.Range("L94").Value = 445
$ws.Range("M94").Value = -160.125
$ws.Range("N94").Value = -1347

# Row data per hunk 17 (old_start 20992)
$ws.Range("H134").Value = 3900.8572
$ws.Range("I134").Value = 1084.4445
$ws.Range("J134").Value = 20799.334
$ws.Range("K134").Value = 3253.3335
$ws.Range("L134").Value = 62398.00199999999
$ws.Range("M134").Value = -718.3335000000002
$ws.Range("N134").Value = -67468.00199999999

$ws = $wb.Worksheets.Item("CRP")
# Row data per hunk 18 (old_start 24297)
$ws.Range("H58").Value = 2880.3125
$ws.Range("I58").Value = 1532.9166
$ws.Range("K58").Value = 1532.9166
$ws.Range("M58").Value = -1329.9166

# Row data per hunk 19 (old_start 27893)
$ws.Range("H132").Value = 2223.4783
$ws.Range("I132").Value = 1959.2858
$ws.Range("K132").Value = 5877.857400000001
$ws.Range("M132").Value = -3347.857400000001

# Row data per hunk 20 (old_start 27991)
$ws.Range("H134").Value = 3016.8572
$ws.Range("I134").Value = 1126
$ws.Range("K134").Value = 3378
$ws.Range("M134").Value = -843

# Row data per hunk 21 (old_start 28092)
$ws.Range("H136").Value = 2880.3125
$ws.Range("I136").Value = 1532.9166
$ws.Range("K136").Value = 4598.7498
$ws.Range("M136").Value = -2048.7498

$ws = $wb.Worksheets.Item("CUL")
# Row data per hunk 22 (old_start 30105)
$ws.Range("H34").Value = 1851.1471
$ws.Range("I34").Value = 186.2
$ws.Range("J34").Value = 2138.2068
$ws.Range("K34").Value = 558.5999999999999
$ws.Range("L34").Value = 6414.6204
$ws.Range("M34").Value = -474.5999999999999
$ws.Range("N34").Value = -6582.6204

# Row data per hunk 23 (old_start 30350)
$ws.Range("H39").Value = 7895.727
$ws.Range("J39").Value = 7895.727
$ws.Range("L39").Value = 23687.181
$ws.Range("N39").Value = -24275.181

# Row data per hunk 24 (old_start 31146)
$ws.Range("H55").Value = 3833.353
$ws.Range("J55").Value = 4535
$ws.Range("L55").Value = 13605
$ws.Range("N55").Value = -13959

# Row data per hunk 25 (old_start 33743)
$ws.Range("H108").Value = 0
$ws.Range("I108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("K108").Value = 0
$ws.Range("L108").Value = 0
$ws.Range("M108").Value = ""
$ws.Range("N108").Value = ""

# Row data per hunk 26 (old_start 33847)
$ws.Range("H110").Value = 6009
$ws.Range("I110").Value = 6009
$ws.Range("K110").Value = 18027
$ws.Range("M110").Value = -13937

# Row data per hunk 27 (old_start 34193)
$ws.Range("H117").Value = 2332.6667
$ws.Range("I117").Value = 0
$ws.Range("K117").Value = 0
$ws.Range("M117").Value = ""

# Row data per hunk 28 (old_start 34343)
$ws.Range("H120").Value = 520
$ws.Range("I120").Value = 520
$ws.Range("K120").Value = 1560
$ws.Range("M120").Value = 3278

$ws = $wb.Worksheets.Item("LTW")
# Row data per hunk 29 (old_start 43414)
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("M23").Value = ""

# Row data per hunk 30 (old_start 45276)
$ws.Range("H61").Value = 4439.4
$ws.Range("I61").Value = 2099.75
$ws.Range("K61").Value = 2099.75
$ws.Range("M61").Value = -1897.75

# Row data per hunk 31 (old_start 46820)
$ws.Range("H93").Value = 1392.5834
$ws.Range("I93").Value = 1437
$ws.Range("K93").Value = 1437
$ws.Range("M93").Value = -189

# Row data per hunk 32 (old_start 47785)
$ws.Range("H113").Value = 4439.4
$ws.Range("I113").Value = 2099.75
$ws.Range("K113").Value = 2099.75
$ws.Range("M113").Value = 70.25

$ws = $wb.Worksheets.Item("WVR")
# Row data per hunk 33 (old_start 51641)
$ws.Range("H51").Value = 23161.666
$ws.Range("J51").Value = 19500
$ws.Range("L51").Value = 19500
$ws.Range("N51").Value = -20520

# Row data per hunk 34 (old_start 52119)
$ws.Range("H61").Value = 7333.3335
$ws.Range("I61").Value = 5500
$ws.Range("K61").Value = 5500
$ws.Range("M61").Value = -5208

# Row data per hunk 35 (old_start 54024)
$ws.Range("H100").Value = 1210.5834
$ws.Range("I100").Value = 1242.7
$ws.Range("J100").Value = 1050
$ws.Range("K100").Value = 2485.4
$ws.Range("L100").Value = 2100
$ws.Range("M100").Value = -1944.4
